# "Cambio limpieza de formato" - clear out the sample quotation data that
# was only there as a template example (CAS / CL-5500 parts, Stock/3 Dias/
# 5 Dias delivery, the "Tarjeta..." descriptions, and the vendor quote
# amounts), leaving the row's structure/styles/formulas intact but blank,
# the same way the author reset the sheet before reusing the template.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Comparativo")

# --- Row 11 ("Registro, Basculas y Servicios" quote line) ---
$ws1.Range("D11").ClearContents()
$ws1.Range("E11").ClearContents()
$ws1.Range("F11").ClearContents()
$ws1.Range("G11").ClearContents()
$ws1.Range("H11").Value = 0
$ws1.Range("N11").ClearContents()

# --- Row 12 ("Teczone" quote line) ---
$ws1.Range("D12").ClearContents()
$ws1.Range("E12").ClearContents()
$ws1.Range("F12").ClearContents()
$ws1.Range("G12").ClearContents()
$ws1.Range("K12").Value = 0
$ws1.Range("N12").ClearContents()

# --- Row 13 ("Bicom" quote line) ---
$ws1.Range("D13").ClearContents()
$ws1.Range("E13").ClearContents()
$ws1.Range("F13").ClearContents()
$ws1.Range("G13").ClearContents()
$ws1.Range("K13").Value = 0
$ws1.Range("N13").ClearContents()

# Match the author's last-saved selection on the Comparativo sheet.
$ws1.Range("N10").Select()
